$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.656247615814209
$ws.Range("B1").Value = 1.643027424812317
$ws.Range("C1").Value = 1.95225715637207
$ws.Range("D1").Value = 3.586067676544189
$ws.Range("E1").Value = 3.845049858093262
